$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - copy format (bold/border/center) from E1, then set value
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data cells F2:F52 - time_taken timestamps (plain/default style, matching column E data cells)
$ws.Range("F2").Value = "2021-10-05 10:50:38.502470"
$ws.Range("F3").Value = "2021-10-05 10:50:38.502481"
$ws.Range("F4").Value = "2021-10-05 10:50:38.502484"
$ws.Range("F5").Value = "2021-10-05 10:50:38.502487"
$ws.Range("F6").Value = "2021-10-05 10:50:38.502490"
$ws.Range("F7").Value = "2021-10-05 10:50:38.502493"
$ws.Range("F8").Value = "2021-10-05 10:50:38.502496"
$ws.Range("F9").Value = "2021-10-05 10:50:38.502498"
$ws.Range("F10").Value = "2021-10-05 10:50:38.502501"
$ws.Range("F11").Value = "2021-10-05 10:50:38.502504"
$ws.Range("F12").Value = "2021-10-05 10:50:38.502506"
$ws.Range("F13").Value = "2021-10-05 10:50:38.502509"
$ws.Range("F14").Value = "2021-10-05 10:50:38.502512"
$ws.Range("F15").Value = "2021-10-05 10:50:38.502514"
$ws.Range("F16").Value = "2021-10-05 10:50:38.502516"
$ws.Range("F17").Value = "2021-10-05 10:50:38.502519"
$ws.Range("F18").Value = "2021-10-05 10:50:38.502522"
$ws.Range("F19").Value = "2021-10-05 10:50:38.502525"
$ws.Range("F20").Value = "2021-10-05 10:50:38.502527"
$ws.Range("F21").Value = "2021-10-05 10:50:38.502530"
$ws.Range("F22").Value = "2021-10-05 10:50:38.502533"
$ws.Range("F23").Value = "2021-10-05 10:50:38.502535"
$ws.Range("F24").Value = "2021-10-05 10:50:38.502538"
$ws.Range("F25").Value = "2021-10-05 10:50:38.502540"
$ws.Range("F26").Value = "2021-10-05 10:50:38.502543"
$ws.Range("F27").Value = "2021-10-05 10:50:38.502546"
$ws.Range("F28").Value = "2021-10-05 10:50:38.502549"
$ws.Range("F29").Value = "2021-10-05 10:50:38.502551"
$ws.Range("F30").Value = "2021-10-05 10:50:38.502554"
$ws.Range("F31").Value = "2021-10-05 10:50:38.502557"
$ws.Range("F32").Value = "2021-10-05 10:50:38.502560"
$ws.Range("F33").Value = "2021-10-05 10:50:38.502562"
$ws.Range("F34").Value = "2021-10-05 10:50:38.502565"
$ws.Range("F35").Value = "2021-10-05 10:50:38.502568"
$ws.Range("F36").Value = "2021-10-05 10:50:38.502571"
$ws.Range("F37").Value = "2021-10-05 10:50:38.502574"
$ws.Range("F38").Value = "2021-10-05 10:50:38.502576"
$ws.Range("F39").Value = "2021-10-05 10:50:38.502579"
$ws.Range("F40").Value = "2021-10-05 10:50:38.502581"
$ws.Range("F41").Value = "2021-10-05 10:50:38.502585"
$ws.Range("F42").Value = "2021-10-05 10:50:38.502590"
$ws.Range("F43").Value = "2021-10-05 10:50:38.502595"
$ws.Range("F44").Value = "2021-10-05 10:50:38.502600"
$ws.Range("F45").Value = "2021-10-05 10:50:38.502605"
$ws.Range("F46").Value = "2021-10-05 10:50:38.502609"
$ws.Range("F47").Value = "2021-10-05 10:50:38.502614"
$ws.Range("F48").Value = "2021-10-05 10:50:38.502617"
$ws.Range("F49").Value = "2021-10-05 10:50:38.502619"
$ws.Range("F50").Value = "2021-10-05 10:50:38.502622"
$ws.Range("F51").Value = "2021-10-05 10:50:38.502625"
$ws.Range("F52").Value = "2021-10-05 10:50:38.502627"
